# Update cost results row 2 values across all year sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("2025", "2030", "2035", "2040", "2045", "2050")

$data = @{
    "2025" = @{
        "A2"=0;  "B2"=1037.265132737054; "C2"=0; "D2"=0;
        "E2"=28926.05393052954; "F2"=0;
        "G2"=8095.925712661834; "H2"=0;
        "I2"=16171.06685703679; "J2"=0; "K2"=0;
        "L2"=48492.22142001599; "M2"=10595.37713982;
        "N2"=7074.779504295464; "O2"=6994.715574734591
    }
    "2030" = @{
        "A2"=0; "B2"=4157.588990853394; "C2"=0; "D2"=0;
        "E2"=45991.90904307188; "F2"=0;
        "G2"=8095.925712661834; "H2"=0;
        "I2"=37079.12819938764; "J2"=0; "K2"=0;
        "L2"=54844.03303316472; "M2"=17449.04999683176;
        "N2"=9029.080117872556; "O2"=9726.534234532202
    }
    "2035" = @{
        "A2"=2754.31755456332; "B2"=6368.910634126893; "C2"=0; "D2"=0;
        "E2"=57457.45307013817; "F2"=0;
        "G2"=8095.925712661834; "H2"=0;
        "I2"=52465.73681402855; "J2"=0; "K2"=0;
        "L2"=54844.03303316472; "M2"=21912.87293902603;
        "N2"=13040.8373108866; "O2"=12862.69800702035
    }
    "2040" = @{
        "A2"=2754.31755456332; "B2"=6368.910634126893; "C2"=0; "D2"=0;
        "E2"=57457.45307013817; "F2"=0;
        "G2"=8095.925712661834; "H2"=0;
        "I2"=52465.73681402855; "J2"=0; "K2"=0;
        "L2"=54844.03303316472; "M2"=21912.87293902603;
        "N2"=13158.56380735076; "O2"=12862.69800702035
    }
    "2045" = @{
        "A2"=5713.151062849596; "B2"=6368.910634126893; "C2"=0; "D2"=0;
        "E2"=57457.45307013817; "F2"=0;
        "G2"=8095.925712661834; "H2"=0;
        "I2"=52465.73681402855; "J2"=0; "K2"=0;
        "L2"=54844.03303316472; "M2"=21912.87293902603;
        "N2"=13608.52715637408; "O2"=14941.40968327155
    }
    "2050" = @{
        "A2"=5713.151062849596; "B2"=6368.910634126893; "C2"=0; "D2"=0;
        "E2"=57457.45307013817; "F2"=0;
        "G2"=8095.925712661834; "H2"=0;
        "I2"=52465.73681402855; "J2"=0; "K2"=0;
        "L2"=54844.03303316472; "M2"=21912.87293902603;
        "N2"=13608.52715637408; "O2"=14941.40968327155
    }
}

foreach ($sheetName in $sheetNames) {
    $sname = [string]$sheetName
    $ws = $wb.Worksheets.Item($sname)
    $rowData = $data[$sname]
    foreach ($cellRef in $rowData.Keys) {
        $ws.Range($cellRef).Value = $rowData[$cellRef]
    }
}
